# Add a new "aluminum polymer" (22uF) capacitor line item to the BOM.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7 (pushes existing row 7 and everything below it
# down by one row, and Excel auto-adjusts the SUMPRODUCT formula ranges).
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new part.
$ws.Cells.Item(7, 1).Value = "22uF"
$ws.Cells.Item(7, 2).Value = 1
$ws.Cells.Item(7, 3).Value = "P16405CT-ND"
$ws.Cells.Item(7, 4).Value = 0.989
$ws.Cells.Item(7, 5).Value = 0.589

# The row-insert operation leaves row 21 ("Board (est.)") column F with the
# number format copied from the old row 20 (the "$...409 locale" format).
# In the target workbook that cell instead uses the plain "$" currency
# format (same as column D/E), so fix it up explicitly.
$ws.Cells.Item(21, 6).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat

# Update the active selection to match the saved view state.
$ws.Range("C18").Select()
